$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.709.25"
$ws.Range("E2").Value = "  +0.25%  "

$ws.Range("D3").Value = "1.920.20"
$ws.Range("E3").Value = "  +1.64%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4936"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.39%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3011"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.96%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06781"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("D10").Value = "1.910.62"
$ws.Range("E10").Value = "  +1.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.27"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07333"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.35%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.236"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.80%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.79%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6755"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.42%  "

$ws.Range("D16").Value = "30.681.87"
$ws.Range("E16").Value = "  +0.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007987"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.24%  "

$ws.Range("E19").Value = "  +0.05%  "

$ws.Range("D20").Value = "2.157.14"
$ws.Range("E20").Value = "  +1.20%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.395"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +11.87%  "

$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "199.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.346"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.687"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.74%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.966"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.476"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.374"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09178"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.083"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05280"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.58%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7453"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.46%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.123"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.99%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.708"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.26%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01856"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.98%  "

$ws.Range("E38").Value = "  +1.99%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9281"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.91%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.091"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.65%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4509"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.98%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +25.66%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.951"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.41%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "107.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.91%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1406"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.36%  "

$ws.Range("E46").Value = "  +0.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.720"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.066"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.75%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.44%  "

$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4060"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.48%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05892"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.60%  "
